$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G10" = 2.45
    "H10" = 3.25
    "I10" = 2.77
    "J10" = 3
    "K10" = 2.07
    "L10" = 3.4
    "O10" = 1.31
    "P10" = 3.15
    "Q10" = 1.93
    "R10" = 1.8
    "T10" = 2.67
    "U10" = 1.72
    "V10" = 2
    "X10" = 12.5
    "Y10" = 9.25
    "Z10" = 26
    "AA10" = 19.5
    "AB10" = 28
    "AD10" = 6.2
    "AE10" = 13.5
    "AF10" = 60
    "AG10" = 450
    "AH10" = 8.5
    "AJ10" = 10
    "AK10" = 32
    "AL10" = 24
    "AM10" = 32
    "AN10" = 4.35
    "AO10" = 13
    "AP10" = 21
    "AQ10" = 55
    "AR10" = 90
    "AS10" = 250
    "AT10" = 2.67
    "AU10" = 7
    "AV10" = 65
    "AW10" = 4.7
    "AX10" = 15.5
    "AY10" = 23
    "AZ10" = 70
    "BA10" = 110
    "G11" = 1.44
    "H11" = 4.2
    "I11" = 6.8
    "J11" = 1.98
    "K11" = 2.25
    "L11" = 6.4
    "N11" = 7.6
    "O11" = 1.28
    "P11" = 3.35
    "Q11" = 1.85
    "R11" = 1.88
    "S11" = 1.4
    "T11" = 2.75
    "X11" = 6.3
    "Y11" = 8.25
    "Z11" = 9.25
    "AA11" = 12.5
    "AB11" = 30
    "AC11" = 7.6
    "AD11" = 8.25
    "AH11" = 16
    "AJ11" = 21
    "AM11" = 75
    "AO11" = 6.8
    "AQ11" = 21
    "AT11" = 2.75
    "AV11" = 90
    "AW11" = 7.8
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
